$d = $word.ActiveDocument

# --- Paragraph text updates (appear 4x, once per "Configuration ..." section) ---
$d.Content.Find.Execute(
    "Operator: Operator:", $false, $false, $false, $false, $false,
    $true, 1, $false,
    "Operator: VA/FK, 17/02/2025 17:31:53", 2)

$d.Content.Find.Execute(
    "Test Configuration: Test Configuration:", $false, $false, $false, $false, $false,
    $true, 1, $false,
    "Test Configuration: Antenna position:in front of harness, DUT Orientation:axis X, Housing connected to the ground plane:no, Configuration of the power return line LV:remotely grounded, -", 2)

$d.Content.Find.Execute(
    "Operating mode: Operating mode:", $false, $false, $false, $false, $false,
    $true, 1, $false,
    "Operating mode: Mode 3, Conclusion  comply", 2)

# --- Results table (first table in the document) ---
$t = $d.Tables.Item(1)

function Set-Cell($table, $row, $col, $text) {
    $rng = $table.Cell($row, $col).Range
    $rng.Text = $text
}

function Set-CellOk($table, $row, $col) {
    $cell = $table.Cell($row, $col)
    $rng = $cell.Range
    $rng.Text = "OK"
    $colored = $d.Range($rng.Start, $rng.Start + 2)
    $colored.Font.Color = 32768
}

# Row 2 : Peak 1.25875
Set-Cell $t 2 6 "33.25"
Set-Cell $t 2 7 "64.0"
Set-Cell $t 2 8 "30.75"
Set-CellOk $t 2 9

# Row 3 : Peak 27.269 -> Q-Peak
Set-Cell $t 3 1 "Q-Peak"
Set-Cell $t 3 6 "22.15"
Set-Cell $t 3 7 "58.0"
Set-Cell $t 3 8 "35.85"
Set-CellOk $t 3 9

# Row 4 : Peak 1.25650 -> Q-Peak
Set-Cell $t 4 1 "Q-Peak"
Set-Cell $t 4 6 "31.39"
Set-Cell $t 4 7 "51.0"
Set-Cell $t 4 8 "19.61"
Set-CellOk $t 4 9

# Row 5 : Peak 27.260
Set-Cell $t 5 6 "16.05"
Set-Cell $t 5 7 "45.0"
Set-Cell $t 5 8 "28.95"
Set-CellOk $t 5 9

# Row 6 : Peak 1.25875
Set-Cell $t 6 6 "33.25"
Set-Cell $t 6 7 "51.0"
Set-Cell $t 6 8 "17.75"
Set-CellOk $t 6 9

# Row 7 : Peak 27.269
Set-Cell $t 7 6 "22.15"
Set-Cell $t 7 7 "45.0"
Set-Cell $t 7 8 "22.85"
Set-CellOk $t 7 9

# Row 8 : CISPR.AVG 1.25425 (already OK / Limite 44.0 unchanged)
Set-Cell $t 8 6 "18.35"
Set-Cell $t 8 8 "25.65"

# Row 9 : CISPR.AVG 27.260 (already OK / Limite 28.0 unchanged)
Set-Cell $t 9 6 "4.57"
Set-Cell $t 9 8 "23.43"
